# Add data for 2021-09-21: the "through September 12" snapshot becomes
# "through September 13" -- the sheet title, the column-B header label, and
# the cumulative counts for any neighborhood/year-column combination that
# saw a new carjacking recorded on 2021-09-13 (and its equivalent date in
# prior years' "through Sept N" columns) are bumped by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item(1)

# Rename the worksheet tab.
$ws.Name = "Through 2021-09-13"

# Update the header label for the current-year column (B1's shared string).
$ws.Range("B1").Value = "September 2021 (through September 13)"

# --- Garfield Park (row 2) ---
$ws.Range("AU2").Value = 1

# --- Humboldt Park (row 4) ---
$ws.Range("K4").Value = 5

# --- Austin (row 5) ---
$ws.Range("B5").Value = 5
$ws.Range("T5").Value = 1
$ws.Range("AC5").Value = 4
$ws.Range("AL5").Value = 2

# --- Roseland (row 6) ---
$ws.Range("T6").Value = 2

# --- Auburn Gresham (row 7) ---
$ws.Range("B7").Value = 6

# --- West Town (row 10) ---
$ws.Range("AL10").Value = 2
$ws.Range("AU10").Value = 4

# --- Little Italy, UIC (row 11) ---
$ws.Range("B11").Value = 2

# --- Chatham (row 13) ---
$ws.Range("B13").Value = 3
$ws.Range("K13").Value = 1

# --- Bucktown (row 15) ---
$ws.Range("B15").Value = 1

# --- Wicker Park (row 19) ---
$ws.Range("B19").Value = 2

# --- Englewood (row 20) ---
$ws.Range("B20").Value = 1

# --- United Center (row 23) ---
$ws.Range("B23").Value = 1
$ws.Range("K23").Value = 2

# --- Ashburn (row 24) ---
$ws.Range("AU24").Value = 3

# --- Logan Square (row 28) ---
$ws.Range("K28").Value = 1

# --- West Loop (row 31) ---
$ws.Range("K31").Value = 4

# --- Lake View (row 33) ---
$ws.Range("K33").Value = 1

# --- Douglas (row 39) ---
$ws.Range("K39").Value = 1

# --- Hegewisch (row 76) ---
$ws.Range("AC76").Value = 1

# --- Printers Row (row 90) ---
$ws.Range("BD90").Value = 1
